$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/postgraduate-scholarships",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/partnership-projects",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/development-grants",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/cre",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/ideas-grants",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/investigator-grants",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/synergy-grants",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/clinical-trials-and-cohort-studies-grants",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/manage-your-grant",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/eucrg",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/nhmrc-e-asia-joint-research-program-grants",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/ideas-grants#Funding%20Partners",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/nhmrc-collaborations-in-health-services-research",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/nhmrc-amed-2024-aspire",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/2024-nhmrc-cihr-ccna",
    "https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/nhmrc-gacd-funding-call"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $urls[$i]
    $ws.Cells.Item($row, 2).Value = "a%"
}
